# Applies the GWDS_explanation.docx edits described in the commit diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Demographic-model sentence: add '(Hap)FLK' to the list of comparable tools
#    'Fdist approaches and Bayescan)' -> 'Fdist approaches, Bayescan and (Hap)FLK)'
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("approaches and Bayescan), GWDS", $true, $false, $false, $false, $false, $true, 1, $false, "approaches, Bayescan and (Hap)FLK), GWDS", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Replace the 'less confounded by assumptions...' sentence at the end of the
#    paragraph with the new 'One major assumption...' text (content moved up
#    from the old following paragraph, and lightly re-worded).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This direct approach is less confounded by assumptions – e.g. no need for assumptions about population hierarchy, equal population sizes, levels of gene flow, etc.", $true, $false, $false, $false, $false, $true, 1, $false, "One major assumption of GWDS is that the vast majority of SNPs are neutral, and that therefore the neutral distribution can be inferred from the overall distribution. On a total of thousands of SNPs, one or two SNPs under selection will not massively influence the overall mean, and hence not affect the rate parameter of the fitted exponential distribution. This assumption is however violated by dense SNP datasets in which many SNPs can represent the same selective sweep. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Heading: 'OPTIONAL: thin data prior to STEP2' -> 'Assumptions underlying GWDS'
#    Also drop italic formatting and the explicit SpaceAfter=0 override on that
#    paragraph (falls back to the 8pt document default).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("OPTIONAL: thin data prior to STEP2", $true, $false, $false, $false, $false, $true, 1, $false, "Assumptions underlying GWDS", 2) | Out-Null
$pHeading = $d.Paragraphs.Item(69)
$pHeading.Range.Font.Italic = 0
$pHeading.Range.Font.ItalicBi = 0
$pHeading.Range.ParagraphFormat.SpaceAfter = 8

# ---------------------------------------------------------------------------
# 4) Former 'One major assumption...' paragraph becomes a new bold sub-heading
#    'Data thinning' (its original sentence content was merged into the
#    paragraph above in step 2).
# ---------------------------------------------------------------------------
$pDataThinning = $d.Paragraphs.Item(71)
$rDT = $pDataThinning.Range
$rDT2 = $d.Range($rDT.Start, $rDT.End - 1)
$rDT2.Text = "Data thinning"
$rDT.Font.Bold = 1
$rDT.Font.BoldBi = 1

# ---------------------------------------------------------------------------
# 5) Rewrite the data-thinning explanation paragraph: the option is no longer
#    offered, drop the OutFLANK/Fst-trimming comparison sentence, and note the
#    option did not improve results.
# ---------------------------------------------------------------------------
$pThinDesc = $d.Paragraphs.Item(72)
$rTD = $pThinDesc.Range
$rTD2 = $d.Range($rTD.Start, $rTD.End - 1)
$rTD2.Text = "GWDS no longer offers the option to infer the neutral distribution from a thinned dataset which contains at maximum 1 SNP per 1MB (or other user defined size). This option was found not to improve the results. "

Write-Output "edits applied"
